$d = $word.ActiveDocument
$paras = $d.Paragraphs
$idx1 = $null
for ($i = 1; $i -le $paras.Count; $i++) {
  $p = $paras.Item($i)
  $t = $p.Range.Text
  if ($t -like "4.Admin wallet*") {
    $idx1 = $i
    break
  }
}
$p1 = $paras.Item($idx1)
$r1 = $p1.Range
$xml1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">4.Admin wallet/add money </w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

$paras2 = $d.Paragraphs
$p2 = $paras2.Item($idx1 + 1)
Write-Host "p2 text before step2=[$($p2.Range.Text)]"
$r2 = $p2.Range
$xml2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>5.Referral Code Integration</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>(**trying)</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>6.wallet recharge by payment gateway</w:t></w:r><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-GB"/></w:rPr><w:t>(**trying)</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)
Write-Host "Final count=$($d.Paragraphs.Count)"
